$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# web_epi rows (2 and 3): schedule changed to Friday-only extracts.
# Clear all day-of-week "X" marks except Friday (column M) on rows 2 and 3.
$ws.Range("H2:L2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H3:L3").ClearContents()
$ws.Range("N3").ClearContents()

# trend_epi row (4): now also runs on Thursday, in addition to Monday.
$ws.Range("L4").Value = "X"

# Move the active selection to the newly-used cell.
$ws.Range("L4").Select()
